$p = $ppt.ActivePresentation

# ---- Slide 1 ----
$s1 = $p.Slides.Item(1)

# Title 1 (shape id=2): move down (reposition only, size unchanged)
$title = $s1.Shapes.Item(1)
$title.Left = 83.36834716796875
$title.Top = 104.65205383300781

# TextBox 2 (shape id=3): move down (reposition only, size unchanged)
$subtitle = $s1.Shapes.Item(2)
$subtitle.Left = 120.0334701538086
$subtitle.Top = 311.91497802734375

# TextBox 3 (shape id=4, "Presented by : Sharon chebet") removed entirely
$s1.Shapes.Item(3).Delete()

# ---- Slide 9 ----
$s9 = $p.Slides.Item(9)
$title9 = $s9.Shapes.Item(1)
$tr = $title9.TextFrame.TextRange
$tr.Delete()
$tr.Text = "Evaluating Model Performance with ROC Curves"
